# Generate Report for Handoff
# Refresh the handoff identifiers/timestamps across the Overview, zh-cn and
# de-de sheets: a new handoff id (8a4d39d5-7a86-4cbc-b43b-540740e36335)
# replaces the old one (717edec0-a356-4009-850f-f89554b614a8), the target
# xlf file hashes are regenerated, and the handoff timestamps advance.

$wb = $excel.ActiveWorkbook

$oldMd  = "717edec0-a356-4009-850f-f89554b614a8.md"
$newMd  = "8a4d39d5-7a86-4cbc-b43b-540740e36335.md"

$oldZhXlf = "717edec0-a356-4009-850f-f89554b614a8.683ab8ea8ad737ef28e32ade502748d8a43da9b4.zh-cn.xlf"
$newZhXlf = "8a4d39d5-7a86-4cbc-b43b-540740e36335.49312c238f1ee8e80aae6a148e2ad73c62c30db2.zh-cn.xlf"

$oldDeXlf = "717edec0-a356-4009-850f-f89554b614a8.683ab8ea8ad737ef28e32ade502748d8a43da9b4.de-de.xlf"
$newDeXlf = "8a4d39d5-7a86-4cbc-b43b-540740e36335.49312c238f1ee8e80aae6a148e2ad73c62c30db2.de-de.xlf"

$newHandoffDate = "2016-03-23 19:13:44"
$newZhDate      = "2016-03-23 19:13:39"
$newDeDate      = "2016-03-23 19:13:44"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = $newHandoffDate
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMd
}

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhDate
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldZhXlf) {
        $hl.TextToDisplay = $newZhXlf
    }
}

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeDate
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldDeXlf) {
        $hl.TextToDisplay = $newDeXlf
    }
}
